# Replace the header row (row 1, columns B:AG) date strings with
# plain numeric offsets 0..31. This turns the "inline string" date
# labels into numbers so the capacity parameters / timestamps can be
# used to retrieve only valid data, per the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startCol = 2   # column B
$endCol = 33    # column AG

for ($col = $startCol; $col -le $endCol; $col++) {
    $value = $col - $startCol
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $value
}
